$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.855.56"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "2.240.81"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.03"
$ws.Range("E5").Value = "  +3.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.28"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.01"
$ws.Range("E10").Value = "  +5.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.71"
$ws.Range("E11").Value = "  +6.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0789"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.55"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "2.587.16"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "2.231.15"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").Value = "40.767.08"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "0.0₃0899"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.22"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "238.92"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("E25").Value = "  +3.26%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.92"
$ws.Range("E28").Value = "  +4.55%  "
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.47"
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.06"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.17"
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.07"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.03"
$ws.Range("E35").Value = "  +4.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0726"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("E38").Value = "  +6.67%  "
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.23"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("E41").Value = "  +4.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.89"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").Value = "2.094.96"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.83"
$ws.Range("E44").Value = "  +8.08%  "
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.13"
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("E47").Value = "  +9.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.81"
$ws.Range("E48").Value = "  -14.89%  "
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("D50").Value = "2.458.97"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("E51").Value = "  +3.08%  "
